# Add team record (Wins/Losses/Ties) columns to the BOS 1994 roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 - copy the formatting already used by the
# existing header row (A1) so the new headers match (bold, centered,
# bordered) without introducing a brand-new style entry.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-47) gets the team's season record repeated across
# the new columns: 54 wins, 61 losses, 0 ties.
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 54
    $ws.Cells.Item($r, 31).Value = 61
    $ws.Cells.Item($r, 32).Value = 0
}
